$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data edit: employee "mhamada" (row 6) refresh/break value corrected from 1 to 0
$ws.Range("E6").Value = 0

# Formatting cleanup: F11:G11 ("שעת_התחלה"/"שעת_סיום" for the last row) drop the
# redundant explicit-fill styling and fall back to the plain bordered/centered
# style already used elsewhere in the sheet (matches cellXfs style index 7).
$rng = $ws.Range("F11:G11")
$rng.HorizontalAlignment = -4108
$rng.Borders.Item(7).LineStyle = 1
$rng.Borders.Item(10).LineStyle = 1

# Leave the selection where the user's last edit was made.
$ws.Range("E6").Select() | Out-Null
